# Apply the 4.3.1.1 "Youth education by gender" update:
#  1. Fix the English title in C1 (drop the stray period after "4.3.1.1").
#  2. Add a new 2021 data column (M), matching the formatting of the
#     existing year columns (picked up from column K, the same way the
#     last "2020" column L was formatted).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Title text correction (C1)
$ws.Range("C1").Value = "4.3.1.1 Youth education by gender"

# 2) Bring column L and the new column M formatting in line with the rest
#    of the yearly data columns (column K).
$ws.Range("K2:K12").Copy()
$ws.Range("L2:L12").PasteSpecial(-4122)
$ws.Range("K2:K12").Copy()
$ws.Range("M2:M12").PasteSpecial(-4122)

# 3) New column M = year 2021 data
$ws.Range("M3").Value = 2021

$ws.Range("M4").Value = 10.8
$ws.Range("M5").Value = 5.2
$ws.Range("M6").Value = 16.2
$ws.Range("M7").Value = 24.2
$ws.Range("M8").Value = 27.6
$ws.Range("M9").Value = 20.9
$ws.Range("M10").Value = 28.5
$ws.Range("M11").Value = 29.7
$ws.Range("M12").Value = 27.5

# 4) Leave the selection where it ended up after the edit
$null = $ws.Range("O2").Select()
